$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("资讯列表") moved from 80% (in progress, due 11.24, with a
# "interface & data adjustments" note) to 100% complete ("--" due date,
# no note) -- i.e. it now matches the same "done" formatting already used
# by row 2. Copy the formatting (fill/number-format) from the completed
# row 2 cells onto row 3, then set the new value/text.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "--"

# The "备注" (note) cell for that row is no longer needed once the item
# is done.
$ws.Range("D3").ClearContents()

# Final selection left on the sheet.
$ws.Range("G8").Select()
